$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues 3-1-24")

# Mark issues #10, #11, #12 (rows 25-27) as Done, matching the formatting
# already used by the other completed rows (strikethrough font).
$ws.Range("A25:L27").Font.Strikethrough = $true

$ws.Range("M25").Value = "Done"
$ws.Range("M26").Value = "Done"
$ws.Range("M27").Value = "Done"
$ws.Range("M25:M27").Font.Strikethrough = $true

# Move the selection to the newly edited cell and scroll the sheet so the
# new rows are in view.
$ws.Activate()
$ws.Range("M27").Select() | Out-Null

# The PivotTables on "Summary of completion" get refreshed as a result of
# the source data edit above; refreshing drops the stray number-format
# override that was left on the totals cells, so clear it explicitly.
$wsSummary = $wb.Worksheets.Item("Summary of completion")
$wsSummary.Range("B3:E3").Style = "Normal"
$wsSummary.Range("B9:E9").Style = "Normal"
$wsSummary.Range("B15:E15").Style = "Normal"
